$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on price cells whose new values would otherwise be
# auto-converted to numbers (losing exact formatting / precision).
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D50').NumberFormat = '@'

$ws.Range('D2').Value = '67.946.98'
$ws.Range('E2').Value = '  +0.64%  '

$ws.Range('D3').Value = '3.740.49'
$ws.Range('E3').Value = '  +0.50%  '

$ws.Range('D4').Value = '0.999'
$ws.Range('E4').Value = '  -0.18%  '

$ws.Range('D5').Value = '591.98'
$ws.Range('E5').Value = '  +0.26%  '

$ws.Range('D6').Value = '166.97'
$ws.Range('E6').Value = '  +1.29%  '

$ws.Range('D7').Value = '3.738.46'
$ws.Range('E7').Value = '  +0.47%  '

$ws.Range('E8').Value = '  -0.03%  '

$ws.Range('D9').Value = '0.521'
$ws.Range('E9').Value = '  +0.80%  '

$ws.Range('E10').Value = '  +0.58%  '

$ws.Range('D11').Value = '6.42'
$ws.Range('E11').Value = '  -0.55%  '

$ws.Range('D12').Value = '0.448'
$ws.Range('E12').Value = '  -0.15%  '

$ws.Range('D13').Value = '0.0000257'
$ws.Range('E13').Value = '  -1.15%  '

$ws.Range('D14').Value = '35.98'
$ws.Range('E14').Value = '  +0.07%  '

$ws.Range('D15').Value = '4.365.31'
$ws.Range('E15').Value = '  +0.43%  '

$ws.Range('D16').Value = '3.743.03'
$ws.Range('E16').Value = '  +0.60%  '

$ws.Range('D17').Value = '67.879.91'
$ws.Range('E17').Value = '  +0.54%  '

$ws.Range('D18').Value = '17.85'
$ws.Range('E18').Value = '  -2.08%  '

$ws.Range('B19').Value = 'Polkadot'
$ws.Range('C19').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D19').Value = '6.98'
$ws.Range('E19').Value = '  -0.31%  '

$ws.Range('B20').Value = 'TRON'
$ws.Range('C20').Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range('D20').Value = '0.112'
$ws.Range('E20').Value = '  +0.58%  '

$ws.Range('D21').Value = '10.61'
$ws.Range('E21').Value = '  -0.47%  '

$ws.Range('D22').Value = '464.21'
$ws.Range('E22').Value = '  -0.41%  '

$ws.Range('D23').Value = '0.695'
$ws.Range('E23').Value = '  -0.43%  '

$ws.Range('B24').Value = 'PEPE'
$ws.Range('C24').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D24').Value = '0.0000147'
$ws.Range('E24').Value = '  +10.47%  '

$ws.Range('B25').Value = 'Litecoin'
$ws.Range('C25').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range('D25').Value = '83.68'
$ws.Range('E25').Value = '  +1.23%  '

$ws.Range('D26').Value = '2.17'
$ws.Range('E26').Value = '  +0.41%  '

$ws.Range('D27').Value = '11.82'
$ws.Range('E27').Value = '  -1.10%  '

$ws.Range('D28').Value = '10.05'
$ws.Range('E28').Value = '  +0.43%  '

$ws.Range('E29').Value = '  +0.05%  '

$ws.Range('E30').Value = '  +0.53%  '

$ws.Range('E31').Value = '  -1.14%  '

$ws.Range('D32').Value = '29.68'
$ws.Range('E32').Value = '  +0.74%  '

$ws.Range('E33').Value = '  -3.31%  '

$ws.Range('D35').Value = '9.09'
$ws.Range('E35').Value = '  +0.98%  '

$ws.Range('D36').Value = '3.694.02'
$ws.Range('E36').Value = '  +0.56%  '

$ws.Range('E37').Value = '  -0.10%  '

$ws.Range('D38').Value = '3.46'
$ws.Range('E38').Value = '  +1.52%  '

$ws.Range('E39').Value = '  +0.26%  '

$ws.Range('D40').Value = '0.992'
$ws.Range('E40').Value = '  +0.30%  '

$ws.Range('D41').Value = '5.75'
$ws.Range('E41').Value = '  +0.52%  '

$ws.Range('D42').Value = '0.999'
$ws.Range('E42').Value = '  -0.18%  '

$ws.Range('D44').Value = '44.15'
$ws.Range('E44').Value = '  +17.26%  '

$ws.Range('D45').Value = '0.299'
$ws.Range('E45').Value = '  -1.60%  '

$ws.Range('D46').Value = '46.75'
$ws.Range('E46').Value = '  +3.41%  '

$ws.Range('D47').Value = '1.89'
$ws.Range('E47').Value = '  -0.43%  '

$ws.Range('D48').Value = '8.39'
$ws.Range('E48').Value = '  -1.40%  '

$ws.Range('D49').Value = '144.66'
$ws.Range('E49').Value = '  +1.00%  '

$ws.Range('D50').Value = '387.73'
$ws.Range('E50').Value = '  +0.91%  '

$ws.Range('D51').Value = '2.760.21'
$ws.Range('E51').Value = '  +3.27%  '
